$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column P (16th column) to make room for the new header
$ws.Columns.Item(16).Insert()

# Set the new header cell value and style (bold black, matching columns C:O)
$newCell = $ws.Range("P1")
$newCell.Value = "Estado actual"

# Copy the formatting (bold, non-red font) from an existing header cell (O1)
# so the new header cell reuses the same style index.
$ws.Range("O1").Copy()
$newCell.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the active selection to match the target state
$ws.Range("P2").Select()
